$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 252-264 ---
$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIRData = @(
    @('2026-01-30', '17:27:56', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:27:58', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:03', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:08', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:13', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:18', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:23', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:28', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:34', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:39', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:44', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:52', '17:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-01-30', '17:28:54', '17:00', 'Bathroom', 'No Motion', 'Inactive')
)
$r = 252
foreach ($row in $wsPIRData) {
    $wsPIR.Cells.Item($r, 1).Value = "'" + $row[0]
    $wsPIR.Cells.Item($r, 2).Value = $row[1]
    $wsPIR.Cells.Item($r, 3).Value = $row[2]
    $wsPIR.Cells.Item($r, 4).Value = $row[3]
    $wsPIR.Cells.Item($r, 5).Value = $row[4]
    $wsPIR.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Humidity sheet: append rows 170-182 ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsHumidityData = @(
    @('2026-01-30', '17:27:56', '17:00', 'Bathroom', '87.4%', 'Active'),
    @('2026-01-30', '17:27:57', '17:00', 'Bathroom', '87.3%', 'Active'),
    @('2026-01-30', '17:27:59', '17:00', 'Bathroom', '87.4%', 'Active'),
    @('2026-01-30', '17:28:04', '17:00', 'Bathroom', '86.4%', 'Active'),
    @('2026-01-30', '17:28:09', '17:00', 'Bathroom', '87.4%', 'Active'),
    @('2026-01-30', '17:28:14', '17:00', 'Bathroom', '87.4%', 'Active'),
    @('2026-01-30', '17:28:19', '17:00', 'Bathroom', '87.4%', 'Active'),
    @('2026-01-30', '17:28:29', '17:00', 'Bathroom', '85.8%', 'Active'),
    @('2026-01-30', '17:28:34', '17:00', 'Bathroom', '87.3%', 'Active'),
    @('2026-01-30', '17:28:39', '17:00', 'Bathroom', '87.3%', 'Active'),
    @('2026-01-30', '17:28:44', '17:00', 'Bathroom', '86.3%', 'Active'),
    @('2026-01-30', '17:28:52', '17:00', 'Bathroom', '87.4%', 'Active'),
    @('2026-01-30', '17:28:54', '17:00', 'Bathroom', '86.4%', 'Active')
)
$r = 170
foreach ($row in $wsHumidityData) {
    $wsHumidity.Cells.Item($r, 1).Value = "'" + $row[0]
    $wsHumidity.Cells.Item($r, 2).Value = $row[1]
    $wsHumidity.Cells.Item($r, 3).Value = $row[2]
    $wsHumidity.Cells.Item($r, 4).Value = $row[3]
    $wsHumidity.Cells.Item($r, 5).Value = "'" + $row[4]
    $wsHumidity.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- mmWave sheet: append rows 49-51 ---
$wsMmwave = $wb.Worksheets.Item("mmWave")
$wsMmwaveData = @(
    @('2026-01-30', '17:28:16', '17:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @('2026-01-30', '17:28:26', '17:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @('2026-01-30', '17:28:51', '17:00', 'Living Room', 'FALL_DETECTED', 'EMERGENCY')
)
$r = 49
foreach ($row in $wsMmwaveData) {
    $wsMmwave.Cells.Item($r, 1).Value = "'" + $row[0]
    $wsMmwave.Cells.Item($r, 2).Value = $row[1]
    $wsMmwave.Cells.Item($r, 3).Value = $row[2]
    $wsMmwave.Cells.Item($r, 4).Value = $row[3]
    $wsMmwave.Cells.Item($r, 5).Value = $row[4]
    $wsMmwave.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
